# ---------------------------------------------------------------------------
# 1) Refresh the cached "today" date placeholder text (4/25/2018 -> 9/3/2018)
#    on the slide master and on every slide layout's Date Placeholder.
# ---------------------------------------------------------------------------
$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Placeholders.Count; $i++) {
        $ph = $shapes.Placeholders.Item($i)
        # ppPlaceholderDate = 16
        if ($ph.PlaceholderFormat.Type -eq 16 -and $ph.HasTextFrame) {
            $tr = $ph.TextFrame.TextRange
            if ($tr.Text -eq "4/25/2018") {
                $tr.Text = "9/3/2018"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Code-sample clean-ups: a couple of adjacent runs that used to be split
#    mid-token get merged back into a single run (same visible text, just
#    fewer <a:r> elements).
# ---------------------------------------------------------------------------
function Merge-Run($shape, [int]$start, [int]$len, [string]$newText) {
    $tr = $shape.TextFrame.TextRange
    $sub = $tr.Characters($start, $len)
    $sub.Text = $newText
}

# Slide 51 & 52: "       string " + "temp = de.path" -> "       string temp = de.path"
$slide51 = $p.Slides.Item(51)
Merge-Run $slide51.Shapes.Item(2) 120 28 "       string temp = de.path"

$slide52 = $p.Slides.Item(52)
Merge-Run $slide52.Shapes.Item(2) 265 28 "       string temp = de.path"

# Slide 57: "if (binary_search (v.begin(), v.end(), 37" + ")) " + "{"
#           -> "if (binary_search (v.begin(), v.end(), 37)) {"
$slide57 = $p.Slides.Item(57)
Merge-Run $slide57.Shapes.Item(2) 183 45 "if (binary_search (v.begin(), v.end(), 37)) {"
